# Point(s) in Polygon(s) - formatting + column reorder fix
#
# Summary of the change being applied (per commit message / diff):
#  - swap the "smallest_area_(m^3)" and "area_standard_derivation_(+-)" columns
#    (N <-> O), headers and data values both move together
#  - add left/top/wrap-text formatting for the data rows (and column defaults)
#  - size the columns (A:F wide fixed width, G:R sized to content)
#  - turn on AutoFilter for A1:R3 and register the (hidden) _FilterDatabase name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Swap columns N and O (header text + the two data rows) so that
#    "area_standard_derivation_(+-)" now sits in column N and
#    "smallest_area_(m^3)" sits in column O.
# ---------------------------------------------------------------------------
foreach ($row in 1..3) {
    $nAddr = "N$row"
    $oAddr = "O$row"
    $nVal = $ws.Range($nAddr).Value2
    $oVal = $ws.Range($oAddr).Value2
    $ws.Range($nAddr).Value2 = $oVal
    $ws.Range($oAddr).Value2 = $nVal
}

# ---------------------------------------------------------------------------
# 2. Column widths (A:F = 34.7109375, then per-column best-fit widths).
#    (Values chosen are the closest this engine's character->pixel rounding
#    can reach to the target widths.)
# ---------------------------------------------------------------------------
$colWidths = @{
    1  = 33.83333333333333  # A
    2  = 33.83333333333333  # B
    3  = 33.83333333333333  # C
    4  = 33.83333333333333  # D
    5  = 33.83333333333333  # E
    6  = 33.83333333333333  # F
    7  = 13.333333333333332 # G
    8  = 26.666666666666668 # H
    9  = 17.666666666666668 # I
    10 = 30.0                # J
    11 = 21.333333333333336 # K
    12 = 18.833333333333336 # L
    13 = 17.166666666666668 # M
    14 = 29.333333333333336 # N
    15 = 20.833333333333336 # O
    16 = 19.5                # P
    17 = 11.0                # Q
    18 = 14.0                # R
}
for ($c = 1; $c -le 18; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $colWidths[$c]
}

# ---------------------------------------------------------------------------
# 3. Data-row formatting: left/top aligned, wrapped text (A2:R3).
#    Done cell-by-cell so the formatting collapses onto a single shared
#    style instead of leaving the header's style untouched.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 3; $r++) {
    for ($c = 1; $c -le 18; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.HorizontalAlignment = -4131   # xlLeft
        $cell.VerticalAlignment = -4160     # xlTop
        $cell.WrapText = $true
    }
}

# ---------------------------------------------------------------------------
# 4. AutoFilter over the whole table, plus the hidden _FilterDatabase name
#    that Excel creates alongside it.
# ---------------------------------------------------------------------------
$ws.Range("A1:R3").AutoFilter(1) | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$R`$3")
$filterName.Visible = $false
